# Journal de bord TPI - add "Partie 4 Realisation" entries (bug popup fix,
# manuel d'installation / utilisation, documentation realisation) and
# update the sheet's current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entryDate = Get-Date -Year 2021 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0
$entryDate = $entryDate.Date

# Row 86: Documentation / Partie 4 Realisation
$ws.Range("B86").Value = "Documentation"
$ws.Range("C86").Value = $entryDate
$ws.Range("D86").Value = 120
$ws.Range("E86").Value = "Partie 4 Realisation"

# Row 87: Resumé du TPI
$ws.Range("B87").Value = "Resumé du TPI"
$ws.Range("C87").Value = $entryDate
$ws.Range("D87").Value = 120

# Row 88: Creation du Manuel d'utilisation
$ws.Range("B88").Value = "Creation du Manuel d'utilisation"
$ws.Range("C88").Value = $entryDate
$ws.Range("D88").Value = 80

# Row 89: Creation du Manuel d'installation
$ws.Range("B89").Value = "Creation du Manuel d'installation"
$ws.Range("C89").Value = $entryDate
$ws.Range("D89").Value = 60

# Fini (new shared string introduced here, referenced by both E88 and E89)
$ws.Range("E88").Value = "Fini"
$ws.Range("E89").Value = "Fini"

# Update the active selection shown when the workbook is reopened
$ws.Range("D90").Select()
